$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 5219.4116
$ws.Range("I6").Value = 5219.4116
$ws.Range("K6").Value = 15658.2348
$ws.Range("M6").Value = -15546.2348
$ws.Range("H12").Value = 349
$ws.Range("I12").Value = 349
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 349
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -179
$ws.Range("N12").Value = $null
$ws.Range("H109").Value = 112000
$ws.Range("J109").Value = 112000
$ws.Range("L109").Value = 112000
$ws.Range("N109").Value = -114774
$ws.Range("H125").Value = 1037.8572
$ws.Range("I125").Value = 487
$ws.Range("J125").Value = 1258.2
$ws.Range("K125").Value = 4383
$ws.Range("L125").Value = 11323.8
$ws.Range("M125").Value = -1923
$ws.Range("N125").Value = -16243.8
$ws.Range("H127").Value = 2821.25
$ws.Range("I127").Value = 2821.25
$ws.Range("K127").Value = 8463.75
$ws.Range("M127").Value = -3503.75
$ws.Range("H132").Value = 2564.4827
$ws.Range("I132").Value = 2752.8845
$ws.Range("J132").Value = 931.6667
$ws.Range("K132").Value = 8258.6535
$ws.Range("L132").Value = 2795.0001
$ws.Range("M132").Value = -5728.6535
$ws.Range("N132").Value = -7855.0001
$ws.Range("H137").Value = 3081.0286
$ws.Range("I137").Value = 2267
$ws.Range("K137").Value = 6801
$ws.Range("M137").Value = -4251
$ws.Range("H138").Value = 2317.7537
$ws.Range("I138").Value = 1308.1034
$ws.Range("J138").Value = 3049.75
$ws.Range("K138").Value = 3924.3102
$ws.Range("L138").Value = 9149.25
$ws.Range("M138").Value = 1215.6898
$ws.Range("N138").Value = -19429.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 2529.4
$ws.Range("I4").Value = 2300
$ws.Range("K4").Value = 2300
$ws.Range("M4").Value = -2184
$ws.Range("H32").Value = 8336626
$ws.Range("I32").Value = 9261809
$ws.Range("K32").Value = 9261809
$ws.Range("M32").Value = -9261522
$ws.Range("H132").Value = 3071.4146
$ws.Range("I132").Value = 1309.96
$ws.Range("K132").Value = 3929.88
$ws.Range("M132").Value = -1399.88

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 6600
$ws.Range("I22").Value = 6600
$ws.Range("K22").Value = 6600
$ws.Range("M22").Value = -6427
$ws.Range("H80").Value = 1901.3334
$ws.Range("I80").Value = 2971
$ws.Range("J80").Value = 965.375
$ws.Range("K80").Value = 2971
$ws.Range("L80").Value = 965.375
$ws.Range("M80").Value = -1973
$ws.Range("N80").Value = -2961.375
$ws.Range("H81").Value = 50887.8
$ws.Range("J81").Value = 50887.8
$ws.Range("L81").Value = 50887.8
$ws.Range("N81").Value = -53009.8
$ws.Range("H83").Value = 1901.3334
$ws.Range("I83").Value = 2971
$ws.Range("J83").Value = 965.375
$ws.Range("K83").Value = 14855
$ws.Range("L83").Value = 4826.875
$ws.Range("M83").Value = -9863
$ws.Range("N83").Value = -14810.875
$ws.Range("H84").Value = 50887.8
$ws.Range("J84").Value = 50887.8
$ws.Range("L84").Value = 152663.4
$ws.Range("N84").Value = -163271.4
$ws.Range("H86").Value = 2607.7144
$ws.Range("I86").Value = 2951.2856
$ws.Range("K86").Value = 2951.2856
$ws.Range("M86").Value = -1828.2856
$ws.Range("H89").Value = 2607.7144
$ws.Range("I89").Value = 2951.2856
$ws.Range("K89").Value = 14756.428
$ws.Range("M89").Value = -9140.428
$ws.Range("H105").Value = 3005.2
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 3005.2
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 3005.2
$ws.Range("M105").Value = $null
$ws.Range("N105").Value = -6499.2
$ws.Range("H134").Value = 40977.383
$ws.Range("I134").Value = 1980.4762
$ws.Range("K134").Value = 5941.4286
$ws.Range("M134").Value = -3406.4286

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 41001.25
$ws.Range("I6").Value = 22000.5
$ws.Range("K6").Value = 22000.5
$ws.Range("M6").Value = -21887.5
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = $null
$ws.Range("N16").Value = $null
$ws.Range("H105").Value = 2766.25
$ws.Range("I105").Value = 2766.25
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 2766.25
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -1019.25
$ws.Range("N105").Value = $null
$ws.Range("H112").Value = 71092.125
$ws.Range("J112").Value = 71092.125
$ws.Range("L112").Value = 71092.125
$ws.Range("N112").Value = -74046.125
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = $null
$ws.Range("N113").Value = $null

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 275
$ws.Range("I63").Value = 275
$ws.Range("K63").Value = 825
$ws.Range("M63").Value = -76
$ws.Range("H66").Value = 275
$ws.Range("I66").Value = 275
$ws.Range("K66").Value = 2475
$ws.Range("M66").Value = 1269
$ws.Range("H104").Value = 5000
$ws.Range("I104").Value = 0
$ws.Range("K104").Value = 0
$ws.Range("M104").Value = $null
$ws.Range("H121").Value = 6044.6
$ws.Range("I121").Value = 2000
$ws.Range("K121").Value = 6000
$ws.Range("M121").Value = -4690
$ws.Range("H136").Value = 6590.357
$ws.Range("I136").Value = 5828.077
$ws.Range("J136").Value = 16500
$ws.Range("K136").Value = 17484.231
$ws.Range("L136").Value = 49500
$ws.Range("M136").Value = -12384.231
$ws.Range("N136").Value = -59700
$ws.Range("H140").Value = 302059.4
$ws.Range("I140").Value = 302059.4
$ws.Range("K140").Value = 906178.2000000001
$ws.Range("M140").Value = -900998.2000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 94031464
$ws.Range("I11").Value = 147653000
$ws.Range("J11").Value = 8237001.5
$ws.Range("K11").Value = 147653000
$ws.Range("L11").Value = 8237001.5
$ws.Range("M11").Value = -147652861
$ws.Range("N11").Value = -8237279.5
$ws.Range("H112").Value = 97950
$ws.Range("J112").Value = 97950
$ws.Range("L112").Value = 97950
$ws.Range("N112").Value = -100166
$ws.Range("H113").Value = 833.625
$ws.Range("I113").Value = 833.625
$ws.Range("K113").Value = 833.625
$ws.Range("M113").Value = 1336.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3929.0356
$ws.Range("I46").Value = 1747.5333
$ws.Range("K46").Value = 1747.5333
$ws.Range("M46").Value = -1559.5333
$ws.Range("H61").Value = 826.13635
$ws.Range("I61").Value = 799.25
$ws.Range("K61").Value = 799.25
$ws.Range("M61").Value = -597.25
$ws.Range("H113").Value = 826.13635
$ws.Range("I113").Value = 799.25
$ws.Range("K113").Value = 799.25
$ws.Range("M113").Value = 1370.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").Value = $null
$ws.Range("H70").Value = 56007
$ws.Range("J70").Value = 61919
$ws.Range("L70").Value = 61919
$ws.Range("N70").Value = -62549
$ws.Range("H73").Value = 56007
$ws.Range("J73").Value = 61919
$ws.Range("L73").Value = 61919
$ws.Range("N73").Value = -64103
$ws.Range("H100").Value = 1662.7858
$ws.Range("I100").Value = 1636.8462
$ws.Range("K100").Value = 3273.6924
$ws.Range("M100").Value = -2732.6924
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").Value = $null
$ws.Range("H132").Value = 3767.8462
$ws.Range("I132").Value = 3180.2727
$ws.Range("K132").Value = 9540.8181
$ws.Range("M132").Value = -7010.8181
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").Value = $null
